# SAE-303 - Update dataset
# Adds a new "Toutes" worksheet summarising the top-5 sports (all regions
# combined) by number of licensed members, appended after "980 - Monaco".

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Toutes"

# Header row (same columns as every other regional sheet)
$newSheet.Range("A1").Value = "région"
$newSheet.Range("B1").Value = "fédération"
$newSheet.Range("C1").Value = "sexe"
$newSheet.Range("D1").Value = "nombre_licencies"
$newSheet.Range("E1").Value = "code_region"
$newSheet.Range("F1").Value = "annee"

# Top-5 sports (all of France, both sexes combined)
$newSheet.Range("A2").Value = "Toutes"
$newSheet.Range("B2").Value = "Fédération Française de Football"
$newSheet.Range("C2").Value = "Tous"
$newSheet.Range("D2").Value = 2215848

$newSheet.Range("A3").Value = "Toutes"
$newSheet.Range("B3").Value = "Fédération Française de Tennis"
$newSheet.Range("C3").Value = "Tous"
$newSheet.Range("D3").Value = 1106989

$newSheet.Range("A4").Value = "Toutes"
$newSheet.Range("B4").Value = "Fédération Française d'Équitation"
$newSheet.Range("C4").Value = "Tous"
$newSheet.Range("D4").Value = 675186

$newSheet.Range("A5").Value = "Toutes"
$newSheet.Range("B5").Value = "Fédération Française de Basketball"
$newSheet.Range("C5").Value = "Tous"
$newSheet.Range("D5").Value = 594408

$newSheet.Range("A6").Value = "Toutes"
$newSheet.Range("B6").Value = "Fédération Française de Handball"
$newSheet.Range("C6").Value = "Tous"
$newSheet.Range("D6").Value = 531864

# Match the bold/centered/bordered header style used on every other sheet by
# copying the formatting (not the values) from an existing header row.
$styleSource = $wb.Worksheets.Item("980 - Monaco")
$styleSource.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths sized to fit the longest value (header or data) in each
# column, mirroring the other sheets in this workbook.
$newSheet.Columns.Item(1).ColumnWidth = 6.71
$newSheet.Columns.Item(2).ColumnWidth = 34.71
$newSheet.Columns.Item(3).ColumnWidth = 4.71
$newSheet.Columns.Item(4).ColumnWidth = 16.71
$newSheet.Columns.Item(5).ColumnWidth = 11.71
$newSheet.Columns.Item(6).ColumnWidth = 5.71

# Restore the originally active sheet/selection.
$originalActive.Activate()
